# Update PLC data 2025-10-13 13:50:35
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7182
$ws.Range("C3").Value = 162600
$ws.Range("C4").Value = 153604
$ws.Range("C7").Value = 5.53
$ws.Range("C8").Value = 64.65000000000001
